$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 156, shifting existing rows 156-176 down to 157-177.
$ws.Rows.Item(156).Insert()

# New row 156 mirrors the (now shifted-down) row 157 data, but with an updated
# Fecha (date) and Volumen.
$ws.Cells.Item(156, 1).Value = 5
$ws.Cells.Item(156, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(156, 3).Value = "Maule"
$ws.Cells.Item(156, 4).Value = 44474
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(157, 4).NumberFormat
$ws.Cells.Item(156, 5).Value = 7
$ws.Cells.Item(156, 6).Value = 100112006
$ws.Cells.Item(156, 7).Value = "Repollo"
$ws.Cells.Item(156, 8).Value = "Crespo record"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 4000
$ws.Cells.Item(156, 11).Value = 500
$ws.Cells.Item(156, 12).Value = 500
$ws.Cells.Item(156, 13).Value = 500
$ws.Cells.Item(156, 14).Value = "$/unidad"
$ws.Cells.Item(156, 15).Value = "Región del Maule"
$ws.Cells.Item(156, 16).Value = 500
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"
